$d = $word.ActiveDocument

# --- Change A -------------------------------------------------------------
# Insert a new sentence about Advertisement chapters (version 2.0.0.0) right
# before "A visual localized notification..." in the intro paragraph.
$d.Content.Find.Execute(
    "to automagically skip past commercials when playing back recordings in Emby.  A visual localized notification",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "to automagically skip past commercials when playing back recordings in Emby.  Additionally, as of version 2.0.0.0, this plugin will recognize Advertisement chapters added by applications such as PlayOn.  A visual localized notification",
    2
) | Out-Null

# --- Change B -------------------------------------------------------------
# Insert a new sentence about Advertisement chapters being skipped when no
# EDL file exists, right before "So," in the "After the ComSkipper plugin..."
# paragraph.
$d.Content.Find.Execute(
    "the commercial areas defined in EDL file will automatically be skipped.  So,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "the commercial areas defined in EDL file will automatically be skipped.  Additionally, if there is no EDL file, but Advertisement chapters are defined in the file, they will be skipped.  So,",
    2
) | Out-Null

# --- Change C -------------------------------------------------------------
# Mention that chapters can also be inserted, not just EDL files created.
$d.Content.Find.Execute(
    "How the EDL files get created is up to the server owner.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "How the EDL files get created, or the chapters get inserted, is up to the server owner.",
    2
) | Out-Null

# --- Change D -------------------------------------------------------------
# Clarify that the real-time watch-while-recording feature only applies
# when EDL files are used.
$d.Content.Find.Execute(
    "It is possible to watch while recording a show and have the commercials skipped.  To allow for this:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "If EDL files are used, it is possible to watch while recording a show and have the commercials skipped.  To allow for this:",
    2
) | Out-Null
